# Updates the cryptos price/volume table (columns D and E, rows 2-51)
# to the latest scraped values, matching the commit's GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.941.00"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.499.79"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'594.67"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'172.75"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.579"
$ws.Range("E8").Value = "  -1.93%  "
$ws.Range("D9").Value = "'0.132"
$ws.Range("E9").Value = "  +3.50%  "
$ws.Range("D10").Value = "'7.14"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "4.093.64"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'29.44"
$ws.Range("E14").Value = "  +5.01%  "
$ws.Range("D15").Value = "66.870.81"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "3.496.19"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'14.31"
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("D20").Value = "'390.59"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'7.93"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "'73.20"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'0.534"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'10.10"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'6.11"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("D31").Value = "'1.42"
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'23.59"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").Value = "'163.13"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "'0.876"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "'4.63"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "2.820.09"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").Value = "'26.94"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").Value = "'0.0728"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "'25.94"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "'0.0300"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").Value = "'339.37"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").Value = "'33.79"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "'6.42"
$ws.Range("E51").Value = "  -0.63%  "
